# Natmi following Dr Hou advice
# Update Vtn-Itgb3 LR-pair sheet: recompute sending/target cluster rows
# across the full 3x3 cluster combination (ECs, FAPs, sCs) with updated
# NATMI specificity metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Vtn"
$ws.Cells.Item(2, 3).Value = "Itgb3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 7.134618
$ws.Cells.Item(2, 8).Value = 21.403854
$ws.Cells.Item(2, 9).Value = 0.0965317920926077
$ws.Cells.Item(2, 10).Value = 0.0965317920926077
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 8.970048
$ws.Cells.Item(2, 14).Value = 26.910144
$ws.Cells.Item(2, 15).Value = 0.487108783009476
$ws.Cells.Item(2, 16).Value = 0.4871087830094759
$ws.Cells.Item(2, 17).Value = 63.997865921664
$ws.Cells.Item(2, 18).Value = 575.9807932949759
$ws.Cells.Item(2, 19).Value = 0.04702148376795389
$ws.Cells.Item(2, 20).Value = 0.04702148376795389

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Vtn"
$ws.Cells.Item(3, 3).Value = "Itgb3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 7.134618
$ws.Cells.Item(3, 8).Value = 21.403854
$ws.Cells.Item(3, 9).Value = 0.0965317920926077
$ws.Cells.Item(3, 10).Value = 0.0965317920926077
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 9.012070666666666
$ws.Cells.Item(3, 14).Value = 27.036212
$ws.Cells.Item(3, 15).Value = 0.489390778604016
$ws.Cells.Item(3, 16).Value = 0.489390778604016
$ws.Cells.Item(3, 17).Value = 64.29768159567199
$ws.Cells.Item(3, 18).Value = 578.679134361048
$ws.Cells.Item(3, 19).Value = 0.04724176889224228
$ws.Cells.Item(3, 20).Value = 0.04724176889224228

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Vtn"
$ws.Cells.Item(4, 3).Value = "Itgb3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 7.134618
$ws.Cells.Item(4, 8).Value = 21.403854
$ws.Cells.Item(4, 9).Value = 0.0965317920926077
$ws.Cells.Item(4, 10).Value = 0.0965317920926077
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.4327576666666667
$ws.Cells.Item(4, 14).Value = 1.298273
$ws.Cells.Item(4, 15).Value = 0.02350043838650813
$ws.Cells.Item(4, 16).Value = 0.02350043838650813
$ws.Cells.Item(4, 17).Value = 3.087560638238
$ws.Cells.Item(4, 18).Value = 27.788045744142
$ws.Cells.Item(4, 19).Value = 0.00226853943241154
$ws.Cells.Item(4, 20).Value = 0.00226853943241154

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Vtn"
$ws.Cells.Item(5, 3).Value = "Itgb3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 17.50798033333334
$ws.Cells.Item(5, 8).Value = 52.52394100000001
$ws.Cells.Item(5, 9).Value = 0.2368839813846793
$ws.Cells.Item(5, 10).Value = 0.2368839813846794
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 8.970048
$ws.Cells.Item(5, 14).Value = 26.910144
$ws.Cells.Item(5, 15).Value = 0.487108783009476
$ws.Cells.Item(5, 16).Value = 0.4871087830094759
$ws.Cells.Item(5, 17).Value = 157.047423973056
$ws.Cells.Item(5, 18).Value = 1413.426815757504
$ws.Cells.Item(5, 19).Value = 0.1153882678867305
$ws.Cells.Item(5, 20).Value = 0.1153882678867305

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Vtn"
$ws.Cells.Item(6, 3).Value = "Itgb3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 17.50798033333334
$ws.Cells.Item(6, 8).Value = 52.52394100000001
$ws.Cells.Item(6, 9).Value = 0.2368839813846793
$ws.Cells.Item(6, 10).Value = 0.2368839813846794
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 9.012070666666666
$ws.Cells.Item(6, 14).Value = 27.036212
$ws.Cells.Item(6, 15).Value = 0.489390778604016
$ws.Cells.Item(6, 16).Value = 0.489390778604016
$ws.Cells.Item(6, 17).Value = 157.7831559946102
$ws.Cells.Item(6, 18).Value = 1420.048403951492
$ws.Cells.Item(6, 19).Value = 0.1159288360886674
$ws.Cells.Item(6, 20).Value = 0.1159288360886675

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Vtn"
$ws.Cells.Item(7, 3).Value = "Itgb3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 17.50798033333334
$ws.Cells.Item(7, 8).Value = 52.52394100000001
$ws.Cells.Item(7, 9).Value = 0.2368839813846793
$ws.Cells.Item(7, 10).Value = 0.2368839813846794
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.4327576666666667
$ws.Cells.Item(7, 14).Value = 1.298273
$ws.Cells.Item(7, 15).Value = 0.02350043838650813
$ws.Cells.Item(7, 16).Value = 0.02350043838650813
$ws.Cells.Item(7, 17).Value = 7.576712717099223
$ws.Cells.Item(7, 18).Value = 68.19041445389301
$ws.Cells.Item(7, 19).Value = 0.005566877409281395
$ws.Cells.Item(7, 20).Value = 0.005566877409281396

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Vtn"
$ws.Cells.Item(8, 3).Value = "Itgb3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 49.26691733333334
$ws.Cells.Item(8, 8).Value = 147.800752
$ws.Cells.Item(8, 9).Value = 0.6665842265227129
$ws.Cells.Item(8, 10).Value = 0.666584226522713
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 8.970048
$ws.Cells.Item(8, 14).Value = 26.910144
$ws.Cells.Item(8, 15).Value = 0.487108783009476
$ws.Cells.Item(8, 16).Value = 0.4871087830094759
$ws.Cells.Item(8, 17).Value = 441.9266132920321
$ws.Cells.Item(8, 18).Value = 3977.339519628288
$ws.Cells.Item(8, 19).Value = 0.3246990313547915
$ws.Cells.Item(8, 20).Value = 0.3246990313547916

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Vtn"
$ws.Cells.Item(9, 3).Value = "Itgb3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 49.26691733333334
$ws.Cells.Item(9, 8).Value = 147.800752
$ws.Cells.Item(9, 9).Value = 0.6665842265227129
$ws.Cells.Item(9, 10).Value = 0.666584226522713
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 9.012070666666666
$ws.Cells.Item(9, 14).Value = 27.036212
$ws.Cells.Item(9, 15).Value = 0.489390778604016
$ws.Cells.Item(9, 16).Value = 0.489390778604016
$ws.Cells.Item(9, 17).Value = 443.9969405368249
$ws.Cells.Item(9, 18).Value = 3995.972464831425
$ws.Cells.Item(9, 19).Value = 0.3262201736231062
$ws.Cells.Item(9, 20).Value = 0.3262201736231063

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Vtn"
$ws.Cells.Item(10, 3).Value = "Itgb3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 49.26691733333334
$ws.Cells.Item(10, 8).Value = 147.800752
$ws.Cells.Item(10, 9).Value = 0.6665842265227129
$ws.Cells.Item(10, 10).Value = 0.666584226522713
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.4327576666666667
$ws.Cells.Item(10, 14).Value = 1.298273
$ws.Cells.Item(10, 15).Value = 0.02350043838650813
$ws.Cells.Item(10, 16).Value = 0.02350043838650813
$ws.Cells.Item(10, 17).Value = 21.32063618903289
$ws.Cells.Item(10, 18).Value = 191.885725701296
$ws.Cells.Item(10, 19).Value = 0.01566502154481519
$ws.Cells.Item(10, 20).Value = 0.01566502154481519
